$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXBots")
$ws.Range("Z5").Value = "Snipper Bot "
Write-Output "done"
